# Append: 2025-10-13 18:32 JST
# Update the "取得日時" (acquired datetime) column for all data rows on the
# first sheet ("ランサーズ") from the previous run timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-13 18:32:28"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
